# "Generate Report for Handoff"
# The b.md file has moved from "Include" (still in sync / handed back)
# to "Ready for handoff" -- a fresh handoff package
# (b.63290e5768f688058c7b37413b0a5c26c308f864.<locale>.xlf) was produced for
# it, with new handoff timestamps. Reflect that across all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is "b.md". Its zh-cn / de-de status columns and
# the "Latest Handoff Date" column need to show the new handoff.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-28-18 22:28:49"

# ---------------------------------------------------------------------
# zh-cn detail sheet: row 3 ("b.md") gets the new handoff file name
# (column D, also a hyperlink) plus the new handoff datetime (column E),
# and its status (column C) becomes "Ready for handoff".
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-18 22:28:46"
foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de detail sheet: same shape of edit as zh-cn, different file/date.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-18 22:28:49"
foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
